$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column score values per diff
$ws.Range("B3").Value = 2.46
$ws.Range("B6").Value = 1.57
$ws.Range("B11").Value = 2.86
$ws.Range("B25").Value = 1.31
$ws.Range("B28").Value = 2.23
$ws.Range("B36").Value = 3.46
$ws.Range("B37").Value = 1.39
$ws.Range("B41").Value = 1.95
$ws.Range("B42").Value = 1.63
$ws.Range("B43").Value = 2.08
$ws.Range("B47").Value = 3.15
$ws.Range("B60").Value = 1.8
$ws.Range("B62").Value = 1.74
$ws.Range("B67").Value = 3
$ws.Range("B76").Value = 2
$ws.Range("B80").Value = 3
$ws.Range("B82").Value = 1.61
$ws.Range("B84").Value = 2.12
$ws.Range("B94").Value = 1.74
$ws.Range("B107").Value = 1.58
$ws.Range("B110").Value = 2.27
$ws.Range("B114").Value = 2
$ws.Range("B122").Value = 1.17
$ws.Range("B131").Value = 3
$ws.Range("B134").Value = 2
$ws.Range("B144").Value = 2.42
$ws.Range("B161").Value = 2
$ws.Range("B170").Value = 3
$ws.Range("B176").Value = 2.27
$ws.Range("B177").Value = 1.45
$ws.Range("B184").Value = 3
$ws.Range("B199").Value = 2.82
$ws.Range("B211").Value = 1
$ws.Range("B217").Value = 3
$ws.Range("B220").Value = 2
$ws.Range("B222").Value = 4
$ws.Range("B226").Value = 2.36
$ws.Range("B227").Value = 4.43
$ws.Range("B228").Value = 4
$ws.Range("B230").Value = 3.14
$ws.Range("B238").Value = 2
$ws.Range("B240").Value = 2
$ws.Range("B244").Value = 3.05
$ws.Range("B253").Value = 3.31
$ws.Range("B257").Value = 3.96
$ws.Range("B260").Value = 3
$ws.Range("B261").Value = 2.48
$ws.Range("B262").Value = 2.95
$ws.Range("B263").Value = 4
$ws.Range("B268").Value = 1.29
$ws.Range("B269").Value = 1
$ws.Range("B272").Value = 3.33
$ws.Range("B276").Value = 2
$ws.Range("B290").Value = 2.67
$ws.Range("B297").Value = 3
$ws.Range("B298").Value = 3.26
$ws.Range("B311").Value = 1
$ws.Range("B314").Value = 1.91
$ws.Range("B317").Value = 1.54
$ws.Range("B321").Value = 4
$ws.Range("B322").Value = 2
$ws.Range("B328").Value = 2
$ws.Range("B332").Value = 2
$ws.Range("B333").Value = 4
$ws.Range("B343").Value = 2
$ws.Range("B367").Value = 1.38
$ws.Range("B368").Value = 2.82
$ws.Range("B373").Value = 3
$ws.Range("B377").Value = 4.34
$ws.Range("B386").Value = 1.61
$ws.Range("B387").Value = 1.7
$ws.Range("B399").Value = 4
$ws.Range("B405").Value = 2.27
$ws.Range("B408").Value = 2
$ws.Range("B415").Value = 1.83
$ws.Range("B423").Value = 1.98
$ws.Range("B434").Value = 2
$ws.Range("B436").Value = 2
$ws.Range("B442").Value = 1.38
$ws.Range("B446").Value = 3
$ws.Range("B449").Value = 2.89
$ws.Range("B454").Value = 3.14
$ws.Range("B457").Value = 2.73
$ws.Range("B461").Value = 3
$ws.Range("B464").Value = 3
$ws.Range("B467").Value = 2.63
$ws.Range("B474").Value = 2
$ws.Range("B475").Value = 2
$ws.Range("B476").Value = 3
$ws.Range("B479").Value = 3.27
$ws.Range("B487").Value = 2.28

# Append new rows 492 and 493
$ws.Range("A492").Value = "Labour"
$ws.Range("B492").Value = 2
$ws.Range("C492").Value = "Negative"
$ws.Range("A493").Value = "Private Job"
$ws.Range("B493").Value = 3.67
$ws.Range("C493").Value = "Positive"
